$d = $word.ActiveDocument
$tbl = $d.Tables.Item(1)

# Row 1 Col 1: "36+25=" -> "47+24="
$tbl.Cell(1, 1).Range.Text = "47+24="
# Row 1 Col 2: "36+6=" -> "80-65="
$tbl.Cell(1, 2).Range.Text = "80-65="
# Row 1 Col 3: "81-63=" -> "20+30="
$tbl.Cell(1, 3).Range.Text = "20+30="
# Row 1 Col 4: "75-56=" -> "90+1="
$tbl.Cell(1, 4).Range.Text = "90+1="
# Row 1 Col 5: "22+69=" -> "34-7="
$tbl.Cell(1, 5).Range.Text = "34-7="
# Row 2 Col 1: "94-46=" -> "49+34="
$tbl.Cell(2, 1).Range.Text = "49+34="
# Row 2 Col 2: "9+51=" -> "15+17="
$tbl.Cell(2, 2).Range.Text = "15+17="
# Row 2 Col 3: "67+9=" -> "8+85="
$tbl.Cell(2, 3).Range.Text = "8+85="
# Row 2 Col 4: "89+1=" -> "84-6="
$tbl.Cell(2, 4).Range.Text = "84-6="
# Row 2 Col 5: "99-70=" -> "59-14="
$tbl.Cell(2, 5).Range.Text = "59-14="
# Row 3 Col 1: "57-55=" -> "24-9="
$tbl.Cell(3, 1).Range.Text = "24-9="
# Row 3 Col 2: "39-11=" -> "18+42="
$tbl.Cell(3, 2).Range.Text = "18+42="
# Row 3 Col 3: "89-9=" -> "53-19="
$tbl.Cell(3, 3).Range.Text = "53-19="
# Row 3 Col 4: "50+12=" -> "1+49="
$tbl.Cell(3, 4).Range.Text = "1+49="
# Row 3 Col 5: "68-19=" -> "65-3="
$tbl.Cell(3, 5).Range.Text = "65-3="
# Row 4 Col 1: "94-49=" -> "8+90="
$tbl.Cell(4, 1).Range.Text = "8+90="
# Row 4 Col 2: "49-9=" -> "36+27="
$tbl.Cell(4, 2).Range.Text = "36+27="
# Row 4 Col 3: "44+31=" -> "8+25="
$tbl.Cell(4, 3).Range.Text = "8+25="
# Row 4 Col 4: "44-14=" -> "0+78="
$tbl.Cell(4, 4).Range.Text = "0+78="
# Row 4 Col 5: "80-56=" -> "20+13="
$tbl.Cell(4, 5).Range.Text = "20+13="
# Row 5 Col 1: "44+48=" -> "30+15="
$tbl.Cell(5, 1).Range.Text = "30+15="
# Row 5 Col 2: "97-71=" -> "34-20="
$tbl.Cell(5, 2).Range.Text = "34-20="
# Row 5 Col 3: "17+10=" -> "28+16="
$tbl.Cell(5, 3).Range.Text = "28+16="
# Row 5 Col 4: "85-82=" -> "2+53="
$tbl.Cell(5, 4).Range.Text = "2+53="
# Row 5 Col 5: "19+22=" -> "1+20="
$tbl.Cell(5, 5).Range.Text = "1+20="
# Row 6 Col 1: "95-3=" -> "62-14="
$tbl.Cell(6, 1).Range.Text = "62-14="
# Row 6 Col 2: "35-13=" -> "29+18="
$tbl.Cell(6, 2).Range.Text = "29+18="
# Row 6 Col 3: "60+14=" -> "82-75="
$tbl.Cell(6, 3).Range.Text = "82-75="
# Row 6 Col 4: "55-6=" -> "93-78="
$tbl.Cell(6, 4).Range.Text = "93-78="
# Row 6 Col 5: "63-55=" -> "46-35="
$tbl.Cell(6, 5).Range.Text = "46-35="
# Row 7 Col 1: "18-11=" -> "45+34="
$tbl.Cell(7, 1).Range.Text = "45+34="
# Row 7 Col 2: "33-11=" -> "57+30="
$tbl.Cell(7, 2).Range.Text = "57+30="
# Row 7 Col 3: "57+16=" -> "86-78="
$tbl.Cell(7, 3).Range.Text = "86-78="
# Row 7 Col 4: "28+64=" -> "43-17="
$tbl.Cell(7, 4).Range.Text = "43-17="
# Row 7 Col 5: "20+12=" -> "53-13="
$tbl.Cell(7, 5).Range.Text = "53-13="
# Row 8 Col 1: "4+36=" -> "65-47="
$tbl.Cell(8, 1).Range.Text = "65-47="
# Row 8 Col 2: "23+28=" -> "30-21="
$tbl.Cell(8, 2).Range.Text = "30-21="
# Row 8 Col 3: "89-37=" -> "4+66="
$tbl.Cell(8, 3).Range.Text = "4+66="
# Row 8 Col 4: "20+70=" -> "10+34="
$tbl.Cell(8, 4).Range.Text = "10+34="
# Row 8 Col 5: "72-52=" -> "30+31="
$tbl.Cell(8, 5).Range.Text = "30+31="
# Row 9 Col 1: "73-20=" -> "19+63="
$tbl.Cell(9, 1).Range.Text = "19+63="
# Row 9 Col 2: "39+14=" -> "87-11="
$tbl.Cell(9, 2).Range.Text = "87-11="
# Row 9 Col 3: "47-16=" -> "9+1="
$tbl.Cell(9, 3).Range.Text = "9+1="
# Row 9 Col 4: "11+86=" -> "27+51="
$tbl.Cell(9, 4).Range.Text = "27+51="
# Row 9 Col 5: "21+67=" -> "82-3="
$tbl.Cell(9, 5).Range.Text = "82-3="
# Row 10 Col 1: "54+11=" -> "25-7="
$tbl.Cell(10, 1).Range.Text = "25-7="
# Row 10 Col 2: "57-41=" -> "9+24="
$tbl.Cell(10, 2).Range.Text = "9+24="
# Row 10 Col 3: "37+26=" -> "59-38="
$tbl.Cell(10, 3).Range.Text = "59-38="
# Row 10 Col 4: "92-64=" -> "17-12="
$tbl.Cell(10, 4).Range.Text = "17-12="
# Row 10 Col 5: "49+46=" -> "43+50="
$tbl.Cell(10, 5).Range.Text = "43+50="
# Row 11 Col 1: "3+7=" -> "31+64="
$tbl.Cell(11, 1).Range.Text = "31+64="
# Row 11 Col 2: "83+2=" -> "62+30="
$tbl.Cell(11, 2).Range.Text = "62+30="
# Row 11 Col 3: "87-49=" -> "51-30="
$tbl.Cell(11, 3).Range.Text = "51-30="
# Row 11 Col 4: "66-33=" -> "85-83="
$tbl.Cell(11, 4).Range.Text = "85-83="
# Row 11 Col 5: "30-20=" -> "92-84="
$tbl.Cell(11, 5).Range.Text = "92-84="
# Row 12 Col 1: "55+39=" -> "59-11="
$tbl.Cell(12, 1).Range.Text = "59-11="
# Row 12 Col 2: "20+64=" -> "33-8="
$tbl.Cell(12, 2).Range.Text = "33-8="
# Row 12 Col 3: "31+29=" -> "28-17="
$tbl.Cell(12, 3).Range.Text = "28-17="
# Row 12 Col 4: "85-8=" -> "12-11="
$tbl.Cell(12, 4).Range.Text = "12-11="
# Row 12 Col 5: "95-29=" -> "62-40="
$tbl.Cell(12, 5).Range.Text = "62-40="
# Row 13 Col 1: "88-37=" -> "36+53="
$tbl.Cell(13, 1).Range.Text = "36+53="
# Row 13 Col 2: "47+48=" -> "39+37="
$tbl.Cell(13, 2).Range.Text = "39+37="
# Row 13 Col 3: "72+13=" -> "1+76="
$tbl.Cell(13, 3).Range.Text = "1+76="
# Row 13 Col 4: "21-20=" -> "85-30="
$tbl.Cell(13, 4).Range.Text = "85-30="
# Row 13 Col 5: "40-6=" -> "10+32="
$tbl.Cell(13, 5).Range.Text = "10+32="
# Row 14 Col 1: "95-13=" -> "98-83="
$tbl.Cell(14, 1).Range.Text = "98-83="
# Row 14 Col 2: "39-7=" -> "61-38="
$tbl.Cell(14, 2).Range.Text = "61-38="
# Row 14 Col 3: "72-12=" -> "83-50="
$tbl.Cell(14, 3).Range.Text = "83-50="
# Row 14 Col 4: "23+38=" -> "48-42="
$tbl.Cell(14, 4).Range.Text = "48-42="
# Row 14 Col 5: "48-12=" -> "69+15="
$tbl.Cell(14, 5).Range.Text = "69+15="
# Row 15 Col 1: "67-58=" -> "49-29="
$tbl.Cell(15, 1).Range.Text = "49-29="
# Row 15 Col 2: "48-27=" -> "10+61="
$tbl.Cell(15, 2).Range.Text = "10+61="
# Row 15 Col 3: "98-65=" -> "74-14="
$tbl.Cell(15, 3).Range.Text = "74-14="
# Row 15 Col 4: "88-25=" -> "29-1="
$tbl.Cell(15, 4).Range.Text = "29-1="
# Row 15 Col 5: "70-56=" -> "30+23="
$tbl.Cell(15, 5).Range.Text = "30+23="
# Row 16 Col 1: "56-49=" -> "96-87="
$tbl.Cell(16, 1).Range.Text = "96-87="
# Row 16 Col 2: "93-7=" -> "4+41="
$tbl.Cell(16, 2).Range.Text = "4+41="
# Row 16 Col 3: "48-32=" -> "94-5="
$tbl.Cell(16, 3).Range.Text = "94-5="
# Row 16 Col 4: "5+17=" -> "76+6="
$tbl.Cell(16, 4).Range.Text = "76+6="
# Row 16 Col 5: "61-45=" -> "64+22="
$tbl.Cell(16, 5).Range.Text = "64+22="
# Row 17 Col 1: "98-30=" -> "99-29="
$tbl.Cell(17, 1).Range.Text = "99-29="
# Row 17 Col 2: "58-9=" -> "10+58="
$tbl.Cell(17, 2).Range.Text = "10+58="
# Row 17 Col 3: "50+6=" -> "37+4="
$tbl.Cell(17, 3).Range.Text = "37+4="
# Row 17 Col 4: "82-4=" -> "14+55="
$tbl.Cell(17, 4).Range.Text = "14+55="
# Row 17 Col 5: "41-35=" -> "70-16="
$tbl.Cell(17, 5).Range.Text = "70-16="
# Row 18 Col 1: "83+11=" -> "56-54="
$tbl.Cell(18, 1).Range.Text = "56-54="
# Row 18 Col 2: "67-30=" -> "86-45="
$tbl.Cell(18, 2).Range.Text = "86-45="
# Row 18 Col 3: "91-85=" -> "4+65="
$tbl.Cell(18, 3).Range.Text = "4+65="
# Row 18 Col 4: "48-24=" -> "83-5="
$tbl.Cell(18, 4).Range.Text = "83-5="
# Row 18 Col 5: "3+88=" -> "66-49="
$tbl.Cell(18, 5).Range.Text = "66-49="
# Row 19 Col 1: "97-39=" -> "94-71="
$tbl.Cell(19, 1).Range.Text = "94-71="
# Row 19 Col 2: "57-44=" -> "66+16="
$tbl.Cell(19, 2).Range.Text = "66+16="
# Row 19 Col 3: "2+42=" -> "8+56="
$tbl.Cell(19, 3).Range.Text = "8+56="
# Row 19 Col 4: "55+1=" -> "60-29="
$tbl.Cell(19, 4).Range.Text = "60-29="
# Row 19 Col 5: "86-39=" -> "7+62="
$tbl.Cell(19, 5).Range.Text = "7+62="
# Row 20 Col 1: "30+57=" -> "60+23="
$tbl.Cell(20, 1).Range.Text = "60+23="
# Row 20 Col 2: "31+33=" -> "14+80="
$tbl.Cell(20, 2).Range.Text = "14+80="
# Row 20 Col 3: "72-22=" -> "79-48="
$tbl.Cell(20, 3).Range.Text = "79-48="
# Row 20 Col 4: "15+29=" -> "75-58="
$tbl.Cell(20, 4).Range.Text = "75-58="
# Row 20 Col 5: "89-23=" -> "39+23="
$tbl.Cell(20, 5).Range.Text = "39+23="
